$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 607.4167
$ws.Range("I33").Value = 619.30304
$ws.Range("J33").Value = 476.66666
$ws.Range("K33").Value = 619.30304
$ws.Range("L33").Value = 476.66666
$ws.Range("M33").Value = -390.30304
$ws.Range("N33").Value = -934.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 99481.44
$ws.Range("I76").Value = 126154
$ws.Range("J76").Value = 3460.2
$ws.Range("K76").Value = 126154
$ws.Range("L76").Value = 3460.2
$ws.Range("M76").Value = -125839
$ws.Range("N76").Value = -4090.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 99481.44
$ws.Range("I79").Value = 126154
$ws.Range("J79").Value = 3460.2
$ws.Range("K79").Value = 126154
$ws.Range("L79").Value = 3460.2
$ws.Range("M79").Value = -125062
$ws.Range("N79").Value = -5644.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1684
$ws.Range("J112").Value = 1889.2
$ws.Range("L112").Value = 5667.6
$ws.Range("N112").Value = -7883.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 738.2
$ws.Range("J127").Value = 1200
$ws.Range("L127").Value = 3600
$ws.Range("N127").Value = -13520

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 955.9259
$ws.Range("J129").Value = 1096
$ws.Range("L129").Value = 3288
$ws.Range("N129").Value = -13288

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 961.63416
$ws.Range("I135").Value = 619.08826
$ws.Range("K135").Value = 5571.79434
$ws.Range("M135").Value = -3036.79434

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 60607868
$ws.Range("I137").Value = 50001508
$ws.Range("J137").Value = 76925340
$ws.Range("K137").Value = 150004524
$ws.Range("L137").Value = 230776020
$ws.Range("M137").Value = -150001974
$ws.Range("N137").Value = -230781120

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2690.0344
$ws.Range("I138").Value = 1153.5
$ws.Range("J138").Value = 3178.932
$ws.Range("K138").Value = 3460.5
$ws.Range("L138").Value = 9536.795999999998
$ws.Range("M138").Value = 1679.5
$ws.Range("N138").Value = -19816.796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 640.2353000000001
$ws.Range("I2").Value = 615.5454999999999
$ws.Range("J2").Value = 685.5
$ws.Range("K2").Value = 615.5454999999999
$ws.Range("L2").Value = 685.5
$ws.Range("M2").Value = -502.5454999999999
$ws.Range("N2").Value = -911.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15867.099
$ws.Range("I32").Value = 6353.6665
$ws.Range("J32").Value = 27079.357
$ws.Range("K32").Value = 6353.6665
$ws.Range("L32").Value = 27079.357
$ws.Range("M32").Value = -6066.6665
$ws.Range("N32").Value = -27653.357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1782
$ws.Range("I45").Value = 1933.5
$ws.Range("J45").Value = 1516.875
$ws.Range("K45").Value = 1933.5
$ws.Range("L45").Value = 1516.875
$ws.Range("M45").Value = -1556.5
$ws.Range("N45").Value = -2270.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 741.9375
$ws.Range("I61").Value = 694.73334
$ws.Range("J61").Value = 1450
$ws.Range("K61").Value = 694.73334
$ws.Range("L61").Value = 1450
$ws.Range("M61").Value = -482.73334
$ws.Range("N61").Value = -1874

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 640.2353000000001
$ws.Range("I116").Value = 615.5454999999999
$ws.Range("J116").Value = 685.5
$ws.Range("K116").Value = 615.5454999999999
$ws.Range("L116").Value = 685.5
$ws.Range("M116").Value = 1678.4545
$ws.Range("N116").Value = -5273.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1255.4783
$ws.Range("I132").Value = 649.2105
$ws.Range("J132").Value = 4135.25
$ws.Range("K132").Value = 1947.6315
$ws.Range("L132").Value = 12405.75
$ws.Range("M132").Value = 582.3685
$ws.Range("N132").Value = -17465.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 741.9375
$ws.Range("I136").Value = 694.73334
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 2084.20002
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = 465.7999799999998
$ws.Range("N136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 640.2353000000001
$ws.Range("I3").Value = 615.5454999999999
$ws.Range("J3").Value = 685.5
$ws.Range("K3").Value = 615.5454999999999
$ws.Range("L3").Value = 685.5
$ws.Range("M3").Value = -501.5454999999999
$ws.Range("N3").Value = -913.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7350.2334
$ws.Range("I94").Value = 575.5769
$ws.Range("J94").Value = 51385.5
$ws.Range("K94").Value = 575.5769
$ws.Range("L94").Value = 51385.5
$ws.Range("M94").Value = -124.5769
$ws.Range("N94").Value = -52287.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2290.2432
$ws.Range("I105").Value = 661.9666999999999
$ws.Range("J105").Value = 9268.571
$ws.Range("K105").Value = 661.9666999999999
$ws.Range("L105").Value = 9268.571
$ws.Range("N105").Value = -12762.571
$ws.Range("M105").Value = 1085.0333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1470.2826
$ws.Range("I134").Value = 1203.6765
$ws.Range("J134").Value = 2225.6667
$ws.Range("K134").Value = 3611.0295
$ws.Range("L134").Value = 6677.000100000001
$ws.Range("M134").Value = -1076.0295
$ws.Range("N134").Value = -11747.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2897.4856
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2897.4856
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2897.4856
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3487.4856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2897.4856
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2897.4856
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2897.4856
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3301.4856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1949.0465
$ws.Range("I58").Value = 487.46875
$ws.Range("J58").Value = 6200.909
$ws.Range("K58").Value = 487.46875
$ws.Range("L58").Value = 6200.909
$ws.Range("M58").Value = -284.46875
$ws.Range("N58").Value = -6606.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1136.7916
$ws.Range("I132").Value = 831.11365
$ws.Range("J132").Value = 4499.25
$ws.Range("K132").Value = 2493.34095
$ws.Range("L132").Value = 13497.75
$ws.Range("M132").Value = 36.65905000000021
$ws.Range("N132").Value = -18557.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1565.909
$ws.Range("I134").Value = 1537.4736
$ws.Range("J134").Value = 1746
$ws.Range("K134").Value = 4612.4208
$ws.Range("L134").Value = 5238
$ws.Range("M134").Value = -2077.4208
$ws.Range("N134").Value = -10308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1949.0465
$ws.Range("I136").Value = 487.46875
$ws.Range("J136").Value = 6200.909
$ws.Range("K136").Value = 1462.40625
$ws.Range("L136").Value = 18602.727
$ws.Range("M136").Value = 1087.59375
$ws.Range("N136").Value = -23702.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5277.5
$ws.Range("I56").Value = 5277.5
$ws.Range("K56").Value = 5277.5
$ws.Range("M56").Value = -4747.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1357.8572
$ws.Range("I113").Value = 382.5
$ws.Range("J113").Value = 1748
$ws.Range("K113").Value = 1147.5
$ws.Range("L113").Value = 5244
$ws.Range("M113").Value = 1022.5
$ws.Range("N113").Value = -9584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1308188.1
$ws.Range("J131").Value = 1163.275
$ws.Range("L131").Value = 3489.825
$ws.Range("N131").Value = -13569.825

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1796.4546
$ws.Range("I132").Value = 1464
$ws.Range("J132").Value = 2683
$ws.Range("K132").Value = 4392
$ws.Range("L132").Value = 8049
$ws.Range("M132").Value = -1862
$ws.Range("N132").Value = -13109

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9925364
$ws.Range("I132").Value = 16031817
$ws.Range("J132").Value = 2376.2917
$ws.Range("K132").Value = 48095451
$ws.Range("L132").Value = 7128.875100000001
$ws.Range("M132").Value = -48092921
$ws.Range("N132").Value = -12188.8751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2133.3333
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 15000
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -18840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1079.8524
$ws.Range("I132").Value = 748.38464
$ws.Range("J132").Value = 1667.4546
$ws.Range("K132").Value = 2245.15392
$ws.Range("L132").Value = 5002.3638
$ws.Range("M132").Value = 284.8460800000003
$ws.Range("N132").Value = -10062.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 930.87695
$ws.Range("I136").Value = 497.51352
$ws.Range("J136").Value = 1503.5358
$ws.Range("K136").Value = 1492.54056
$ws.Range("L136").Value = 4510.607400000001
$ws.Range("M136").Value = 1057.45944
$ws.Range("N136").Value = -9610.607400000001

